# Region XI_LMS.xlsx edit:
# Insert a new column at the start (A) that holds a row index,
# relabel/uppercase several header cells, and shift all existing
# data one column to the right (handled automatically by Insert()).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 25

# Insert a new column before column A; this shifts all existing
# columns (A..Z) right by one (to B..AA), carrying data validation
# along with them (but not per-cell formats for the brand-new column).
$ws.Columns.Item(1).Insert()

# The newly inserted column A comes in with default formatting; copy
# the formats from column B (which now holds what used to be in A) so
# that the header/data styling lines up with the rest of the row.
$ws.Range("B1:B$lastRow").Copy()
$ws.Range("A1:A$lastRow").PasteSpecial(-4122)  # xlPasteFormats

# Set the new column A width
$ws.Columns.Item(1).ColumnWidth = 23

# New header for column A
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"

# Update (uppercase) the header labels that moved into B1:Z1
$ws.Range("B1").Value = "CATEGORY"
$ws.Range("C1").Value = "REGION"
$ws.Range("D1").Value = "DIVISION"
$ws.Range("E1").Value = "SCHOOL ID"
$ws.Range("F1").Value = "SCHOOL NAME"
$ws.Range("G1").Value = "MUNICIPALITY"
$ws.Range("H1").Value = "LD"
$ws.Range("I1").Value = "NO. OF SITES"
$ws.Range("J1").Value = "NO. OF CL"
$ws.Range("K1").Value = "SCOPE OF WORK"
$ws.Range("L1").Value = "PROJECT ALLOCATION"
$ws.Range("M1").Value = "CONTRACT AMOUNT"
$ws.Range("N1").Value = "STATUS"
$ws.Range("O1").Value = "PERCENTAGE OF COMPLETION"
$ws.Range("P1").Value = " TARGET COMPLETION DATE "
$ws.Range("Q1").Value = "ACTUAL DATE OF COMPLETION"
$ws.Range("R1").Value = "PROJECT ID"
$ws.Range("S1").Value = "CONTRACT ID"
$ws.Range("T1").Value = "ISSUANCE OF INVITATION TO BID"
$ws.Range("U1").Value = "PRE-SUBMISSION CONFERENCE"
$ws.Range("V1").Value = "BID OPENING"
$ws.Range("W1").Value = "ISSUANCE OF RESOLUTION TO AWARD"
$ws.Range("X1").Value = "ISSUANCE OF NOTICE TO PROCEED"
$ws.Range("Y1").Value = "NAME OF CONTRACTOR"
$ws.Range("Z1").Value = "OTHER REMARKS"
# AA1 ("Status as of July 11, 2025") is unchanged by the insert.

# Fill in the new INDEX column values for each data row (2-25)
$indexValues = @{
    2  = 91
    3  = 192
    4  = 193
    5  = 194
    6  = 195
    7  = 196
    8  = 197
    9  = 198
    10 = 370
    11 = 371
    12 = 372
    13 = 373
    14 = 374
    15 = 375
    16 = 529
    17 = 530
    18 = 531
    19 = 532
    20 = 533
    21 = 534
    22 = 535
    23 = 536
    24 = 537
    25 = 538
}

foreach ($row in $indexValues.Keys) {
    $ws.Cells.Item($row, 1).Value = $indexValues[$row]
}
